# Daily attendance processing - 2026-01-05 12:57:02
#
# The "Recorded By" column (G) lists the users/processes that recorded each
# attendance session as a comma-separated string, e.g.
#   "backup@backdoor.com, System, system"
# This pass normalizes the ordering of that list so that any "System"
# (case-insensitive) entries are moved to the front, while the relative
# order of all entries is otherwise preserved, e.g.
#   "System, system, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ",\s*"

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    # Only rewrite if there actually is a "System" entry that isn't already
    # leading every other entry (keeps untouched rows untouched).
    if ($systemParts.Length -gt 0) {
        $newVal = ($systemParts + $otherParts) -join ", "
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
